$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Column D = Price (must stay plain text, not get auto-converted to a number)
# Column E = Volume(1h) (already plain text; values contain spaces/% so they stay text)

Set-TextValue "D2" "40.815.59"

Set-TextValue "D3" "2.385.95"
$ws.Range("E3").Value = "  -3.72%  "

$ws.Range("E4").Value = "  -0.05%  "

Set-TextValue "D5" "314.03"
$ws.Range("E5").Value = "  -1.64%  "

Set-TextValue "D6" "88.18"
$ws.Range("E6").Value = "  -5.46%  "

$ws.Range("E7").Value = "  -4.33%  "

$ws.Range("E8").Value = "  +0.09%  "

Set-TextValue "D9" "0.493"
$ws.Range("E9").Value = "  -4.87%  "

Set-TextValue "D10" "0.0824"
$ws.Range("E10").Value = "  -5.23%  "

Set-TextValue "D11" "31.18"
$ws.Range("E11").Value = "  -6.38%  "

$ws.Range("E12").Value = "  -1.71%  "

Set-TextValue "D13" "2.752.82"
$ws.Range("E13").Value = "  -3.78%  "

$ws.Range("E14").Value = "  -4.14%  "

Set-TextValue "D15" "15.16"
$ws.Range("E15").Value = "  -4.01%  "

Set-TextValue "D16" "2.391.22"
$ws.Range("E16").Value = "  -3.46%  "

Set-TextValue "D17" "0.761"
$ws.Range("E17").Value = "  -4.06%  "

Set-TextValue "D18" "40.708.61"
$ws.Range("E18").Value = "  -2.47%  "

Set-TextValue "D19" "0.0₃0911"
$ws.Range("E19").Value = "  -4.26%  "

$ws.Range("E20").Value = "  -4.73%  "

Set-TextValue "D21" "69.30"
$ws.Range("E21").Value = "  -2.81%  "

Set-TextValue "D22" "10.82"
$ws.Range("E22").Value = "  -4.66%  "

Set-TextValue "D23" "233.67"
$ws.Range("E23").Value = "  -2.62%  "

Set-TextValue "D24" "2.65"
$ws.Range("E24").Value = "  -3.50%  "

$ws.Range("E25").Value = "  +0.05%  "

Set-TextValue "D26" "1.82"
$ws.Range("E26").Value = "  -6.30%  "

Set-TextValue "D27" "23.66"

Set-TextValue "D28" "2.18"
$ws.Range("E28").Value = "  -3.84%  "

Set-TextValue "D29" "9.41"
$ws.Range("E29").Value = "  -4.36%  "

Set-TextValue "D30" "34.00"
$ws.Range("E30").Value = "  -6.16%  "

Set-TextValue "D31" "155.83"
$ws.Range("E31").Value = "  -1.54%  "

$ws.Range("E32").Value = "  +0.11%  "

Set-TextValue "D33" "5.24"
$ws.Range("E33").Value = "  -5.25%  "

Set-TextValue "D34" "0.0732"
$ws.Range("E34").Value = "  -4.71%  "

Set-TextValue "D35" "2.41"
$ws.Range("E35").Value = "  -6.70%  "

$ws.Range("E36").Value = "  -2.04%  "

Set-TextValue "D37" "2.82"
$ws.Range("E37").Value = "  -4.13%  "

Set-TextValue "D38" "16.02"
$ws.Range("E38").Value = "  -8.25%  "

Set-TextValue "D39" "0.0999"
$ws.Range("E39").Value = "  -3.49%  "

$ws.Range("E40").Value = "  -7.79%  "

$ws.Range("E41").Value = "  -5.95%  "

$ws.Range("E42").Value = "  -7.93%  "

Set-TextValue "D43" "1.959.29"
$ws.Range("E43").Value = "  -1.70%  "

$ws.Range("E44").Value = "  -5.04%  "

Set-TextValue "D45" "17.73"
$ws.Range("E45").Value = "  -6.84%  "

Set-TextValue "D46" "2.80"
$ws.Range("E46").Value = "  -6.71%  "

Set-TextValue "D47" "9.35"
$ws.Range("E47").Value = "  -1.10%  "

Set-TextValue "D48" "2.615.88"
$ws.Range("E48").Value = "  -3.75%  "

Set-TextValue "D49" "93.76"
$ws.Range("E49").Value = "  -3.97%  "

Set-TextValue "D50" "72.87"
$ws.Range("E50").Value = "  -2.03%  "

Set-TextValue "D51" "50.87"
$ws.Range("E51").Value = "  -3.69%  "
